$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record at row 14 (pushes existing rows 14-55 down to 15-56).
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value  = 4
$ws.Cells.Item(14, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(14, 3).Value  = "Los Lagos"
$ws.Cells.Item(14, 4).Value  = 44868
$ws.Cells.Item(14, 5).Value  = 10
$ws.Cells.Item(14, 6).Value  = 300000000
$ws.Cells.Item(14, 7).Value  = "Espárragos"
$ws.Cells.Item(14, 8).Value  = "Sin especificar"
$ws.Cells.Item(14, 9).Value  = "Primera"
$ws.Cells.Item(14, 10).Value = 600
$ws.Cells.Item(14, 11).Value = 1500
$ws.Cells.Item(14, 12).Value = 1700
$ws.Cells.Item(14, 13).Value = 1600
$ws.Cells.Item(14, 14).Value = "$/kilo"
$ws.Cells.Item(14, 15).Value = "Provincia de Linares"
$ws.Cells.Item(14, 16).Value = 1600
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"
